$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns I1 ("I0") and J1 ("IF"), copying the format of the
# existing header cell H1 so they pick up the same bold/centered/bordered
# style used by the rest of row 1 (reuses the existing style index instead
# of creating a new one).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

# Add the corresponding data values in row 2
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
